$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.345.94"
$ws.Range("E2").Value = "  +0.39%  "

# Row 3
$ws.Range("D3").Value = "3.674.29"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "'686.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6
$ws.Range("D6").Value = "'159.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.492"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9
$ws.Range("E9").Value = "  -1.43%  "

# Row 10
$ws.Range("E10").Value = "  -2.17%  "

# Row 11
$ws.Range("E11").Value = "  -3.57%  "

# Row 13
$ws.Range("D13").Value = "4.299.37"
$ws.Range("E13").Value = "  +0.04%  "

# Row 14
$ws.Range("D14").Value = "'32.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.70%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.682.08"
$ws.Range("E15").Value = "  +0.08%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "69.372.44"
$ws.Range("E16").Value = "  +0.30%  "

# Row 17
$ws.Range("E17").Value = "  +2.03%  "

# Row 18
$ws.Range("D18").Value = "'15.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.24%  "

# Row 19
$ws.Range("D19").Value = "'6.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.90%  "

# Row 20
$ws.Range("D20").Value = "'469.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.75%  "

# Row 21
$ws.Range("D21").Value = "'9.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "

# Row 22
$ws.Range("D22").Value = "'0.648"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.40%  "

# Row 23
$ws.Range("D23").Value = "'79.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").Value = "3.826.24"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").Value = "'0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.01%  "

# Row 27
$ws.Range("D27").Value = "'10.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.42%  "

# Row 28
$ws.Range("D28").Value = "'9.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.68%  "

# Row 29
$ws.Range("E29").Value = "  -1.24%  "

# Row 30
$ws.Range("E30").Value = "  -5.55%  "

# Row 31
$ws.Range("D31").Value = "'6.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.50%  "

# Row 32
$ws.Range("E32").Value = "  -5.75%  "

# Row 33
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "

# Row 34
$ws.Range("D34").Value = "'26.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.36%  "

# Row 35
$ws.Range("D35").Value = "3.651.87"
$ws.Range("E35").Value = "  +0.30%  "

# Row 36
$ws.Range("E36").Value = "  -1.73%  "

# Row 37
$ws.Range("D37").Value = "'8.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.94%  "

# Row 38
$ws.Range("D38").Value = "'6.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.22%  "

# Row 40
$ws.Range("D40").Value = "'2.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.38%  "

# Row 41
$ws.Range("E41").Value = "  -4.97%  "

# Row 42
$ws.Range("E42").Value = "  +0.11%  "

# Row 43
$ws.Range("D43").Value = "'0.940"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "

# Row 44
$ws.Range("D44").Value = "'166.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.09%  "

# Row 45
$ws.Range("D45").Value = "'47.48"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "  +1.49%  "

# Row 47
$ws.Range("D47").Value = "'2.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.86%  "

# Row 48
$ws.Range("E48").Value = "  +5.59%  "

# Row 49
$ws.Range("E49").Value = "  -0.22%  "

# Row 50
$ws.Range("D50").Value = "'27.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.86%  "

# Row 51
$ws.Range("D51").Value = "'7.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.52%  "
